# Update Stückliste Elektrobox Rahmen
# Applies the changes described in the commit diff:
#  - Row 2 (Arduino Nano) gets a real hyperlink on its Link(URL) cell
#  - Row 3's RS-Online connector is swapped for the new part (0471187 / 471-187),
#    with an updated price and a real hyperlink, plus the "pasted from web" style
#    (border + wrap + indent) carried over onto the article-number cell
#  - Rows 6 and 7 (H-Brücke, Taster) also get real hyperlinks on their Link(URL) cells
#  - Selection moves to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: Arduino Nano - add working hyperlink to the existing reichelt URL
# ---------------------------------------------------------------------------
$b2Text = $ws.Range("B2").Value()
$ws.Hyperlinks.Add($ws.Range("B2"), $b2Text) | Out-Null

# ---------------------------------------------------------------------------
# Row 3: Spannungsverteiler - new RS-Online part number/link + new price
# ---------------------------------------------------------------------------
$newRsUrl = "https://at.rs-online.com/web/p/idc-steckverbinder/0471187"
$ws.Range("B3").Value = $newRsUrl
$ws.Range("C3").Value = 5.29
$ws.Range("E3").Value = "471-187"

$ws.Hyperlinks.Add($ws.Range("B3"), $newRsUrl) | Out-Null

# ---------------------------------------------------------------------------
# Row 6: H-Brücke - add working hyperlink to the existing AZDelivery URL
# ---------------------------------------------------------------------------
$b6Text = $ws.Range("B6").Value()
$ws.Hyperlinks.Add($ws.Range("B6"), $b6Text) | Out-Null

# ---------------------------------------------------------------------------
# Row 7: Taster - add working hyperlink to the existing Conrad URL, keeping
# the cached display text (matches the author's original "display" attribute)
# ---------------------------------------------------------------------------
$b7Text = $ws.Range("B7").Value()
$hl7 = $ws.Hyperlinks.Add($ws.Range("B7"), $b7Text)
$hl7.TextToDisplay = $b7Text

# ---------------------------------------------------------------------------
# Re-apply the existing "Link" cell style (already used by B4) to every cell
# that just received a hyperlink, so they pick up the underline/blue look
# instead of a freshly minted style.
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# E3 ("471-187") now carries the formatting that was pasted in from the
# RS-Online website: a grey medium right border plus wrapped / indented /
# vertically centred left-aligned text.
# ---------------------------------------------------------------------------
$e3 = $ws.Range("E3")
$rightBorder = $e3.Borders.Item(10)
$rightBorder.Weight = -4138
$rightBorder.Color = 12500670
$e3.HorizontalAlignment = -4131
$e3.VerticalAlignment = -4108
$e3.WrapText = $true
$e3.IndentLevel = 1

# ---------------------------------------------------------------------------
# Move the active selection to B3, matching the saved sheet view.
# ---------------------------------------------------------------------------
$ws.Range("B3").Select() | Out-Null
